$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3966.4
$ws.Range("I32").Value = 3939
$ws.Range("J32").Value = 3973.25
$ws.Range("K32").Value = 3939
$ws.Range("L32").Value = 3973.25
$ws.Range("M32").Value = -3613
$ws.Range("N32").Value = -4625.25
$ws.Range("H33").Value = 384.36365
$ws.Range("I33").Value = 392.8
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 392.8
$ws.Range("L33").Value = 300
$ws.Range("M33").Value = -163.8
$ws.Range("N33").Value = -758
$ws.Range("H98").Value = 1950
$ws.Range("I98").Value = 1950
$ws.Range("K98").Value = 1950
$ws.Range("M98").Value = -452
$ws.Range("H122").Value = 1950
$ws.Range("I122").Value = 1950
$ws.Range("K122").Value = 5850
$ws.Range("M122").Value = -3400
$ws.Range("H129").Value = 1720.7858
$ws.Range("J129").Value = 2842
$ws.Range("L129").Value = 8526
$ws.Range("N129").Value = -18526
$ws.Range("H137").Value = 1340.75
$ws.Range("I137").Value = 1128.7142
$ws.Range("K137").Value = 3386.1426
$ws.Range("M137").Value = -836.1425999999997
$ws.Range("H138").Value = 5711.353
$ws.Range("I138").Value = 5499.5
$ws.Range("J138").Value = 5739.6
$ws.Range("K138").Value = 16498.5
$ws.Range("L138").Value = 17218.8
$ws.Range("M138").Value = -11358.5
$ws.Range("N138").Value = -27498.8

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 2504
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 2504
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 2504
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = -2792
$ws.Range("H13").Value = 1000000000
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = ""
$ws.Range("H32").Value = 4816622
$ws.Range("I32").Value = 5004697
$ws.Range("J32").Value = 3500100
$ws.Range("K32").Value = 5004697
$ws.Range("L32").Value = 3500100
$ws.Range("M32").Value = -5004410
$ws.Range("N32").Value = -3500674
$ws.Range("H61").Value = 2255.5557
$ws.Range("I61").Value = 2357.5715
$ws.Range("K61").Value = 2357.5715
$ws.Range("M61").Value = -2145.5715
$ws.Range("H74").Value = 2123.75
$ws.Range("I74").Value = 2123.75
$ws.Range("K74").Value = 2123.75
$ws.Range("M74").Value = -1249.75
$ws.Range("H76").Value = 40000
$ws.Range("J76").Value = 40000
$ws.Range("L76").Value = 40000
$ws.Range("N76").Value = -40676
$ws.Range("H77").Value = 2123.75
$ws.Range("I77").Value = 2123.75
$ws.Range("K77").Value = 10618.75
$ws.Range("M77").Value = -6250.75
$ws.Range("H79").Value = 40000
$ws.Range("J79").Value = 40000
$ws.Range("L79").Value = 40000
$ws.Range("N79").Value = -42340
$ws.Range("H97").Value = 1270.4286
$ws.Range("I97").Value = 1323.5
$ws.Range("K97").Value = 1323.5
$ws.Range("M97").Value = -827.5
$ws.Range("H122").Value = 15173.625
$ws.Range("J122").Value = 3059.5
$ws.Range("L122").Value = 9178.5
$ws.Range("N122").Value = -14078.5
$ws.Range("H132").Value = 865.2
$ws.Range("I132").Value = 850.2222
$ws.Range("K132").Value = 2550.6666
$ws.Range("M132").Value = -20.66660000000002
$ws.Range("H136").Value = 2255.5557
$ws.Range("I136").Value = 2357.5715
$ws.Range("K136").Value = 7072.7145
$ws.Range("M136").Value = -4522.7145

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9001.75
$ws.Range("I20").Value = 8669
$ws.Range("K20").Value = 8669
$ws.Range("M20").Value = -8422
$ws.Range("H105").Value = 3000
$ws.Range("I105").Value = 3000
$ws.Range("K105").Value = 3000
$ws.Range("M105").Value = -1253
$ws.Range("H112").Value = 48000
$ws.Range("J112").Value = 48000
$ws.Range("L112").Value = 48000
$ws.Range("N112").Value = -50954
$ws.Range("H134").Value = 4532.8335
$ws.Range("I134").Value = 4532.8335
$ws.Range("K134").Value = 13598.5005
$ws.Range("M134").Value = -11063.5005

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2328.0967
$ws.Range("I31").Value = 1870.6
$ws.Range("J31").Value = 2545.9524
$ws.Range("K31").Value = 1870.6
$ws.Range("L31").Value = 2545.9524
$ws.Range("M31").Value = -1575.6
$ws.Range("N31").Value = -3135.9524
$ws.Range("H34").Value = 2328.0967
$ws.Range("I34").Value = 1870.6
$ws.Range("J34").Value = 2545.9524
$ws.Range("K34").Value = 1870.6
$ws.Range("L34").Value = 2545.9524
$ws.Range("M34").Value = -1668.6
$ws.Range("N34").Value = -2949.9524
$ws.Range("H58").Value = 3124.5557
$ws.Range("I58").Value = 2543
$ws.Range("K58").Value = 2543
$ws.Range("M58").Value = -2340
$ws.Range("H105").Value = 2888.0908
$ws.Range("J105").Value = 3406.5
$ws.Range("L105").Value = 3406.5
$ws.Range("N105").Value = -6900.5
$ws.Range("H125").Value = 97497.836
$ws.Range("J125").Value = 97497.836
$ws.Range("L125").Value = 97497.836
$ws.Range("N125").Value = -102417.836
$ws.Range("H134").Value = 2484.4443
$ws.Range("I134").Value = 2092
$ws.Range("J134").Value = 2975
$ws.Range("K134").Value = 6276
$ws.Range("L134").Value = 8925
$ws.Range("M134").Value = -3741
$ws.Range("N134").Value = -13995
$ws.Range("H136").Value = 3124.5557
$ws.Range("I136").Value = 2543
$ws.Range("K136").Value = 7629
$ws.Range("M136").Value = -5079

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 214.57143
$ws.Range("I12").Value = 174.66667
$ws.Range("J12").Value = 244.5
$ws.Range("K12").Value = 524.00001
$ws.Range("L12").Value = 733.5
$ws.Range("M12").Value = -351.00001
$ws.Range("N12").Value = -1079.5
$ws.Range("H33").Value = 1000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 6000
$ws.Range("M33").Value = ""
$ws.Range("N33").Value = -6566
$ws.Range("H68").Value = 3676.4666
$ws.Range("J68").Value = 3676.4666
$ws.Range("L68").Value = 11029.3998
$ws.Range("N68").Value = -12651.3998
$ws.Range("H71").Value = 3676.4666
$ws.Range("J71").Value = 3676.4666
$ws.Range("L71").Value = 33088.1994
$ws.Range("N71").Value = -41200.1994
$ws.Range("H107").Value = 1723.75
$ws.Range("I107").Value = 1665.1666
$ws.Range("J107").Value = 1899.5
$ws.Range("K107").Value = 4995.4998
$ws.Range("L107").Value = 5698.5
$ws.Range("M107").Value = -3075.4998
$ws.Range("N107").Value = -9538.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 2112600
$ws.Range("I7").Value = 50000
$ws.Range("J7").Value = 2628250
$ws.Range("K7").Value = 50000
$ws.Range("L7").Value = 2628250
$ws.Range("M7").Value = -49888
$ws.Range("N7").Value = -2628474
$ws.Range("H8").Value = 2112600
$ws.Range("I8").Value = 50000
$ws.Range("J8").Value = 2628250
$ws.Range("K8").Value = 50000
$ws.Range("L8").Value = 2628250
$ws.Range("M8").Value = -49861
$ws.Range("N8").Value = -2628528
$ws.Range("H132").Value = 4315.6313
$ws.Range("I132").Value = 4749.875
$ws.Range("J132").Value = 1999.6666
$ws.Range("K132").Value = 14249.625
$ws.Range("L132").Value = 5998.9998
$ws.Range("M132").Value = -11719.625
$ws.Range("N132").Value = -11058.9998

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7771.8438
$ws.Range("I122").Value = 9637.666999999999
$ws.Range("K122").Value = 28913.001
$ws.Range("M122").Value = -26463.001
$ws.Range("H125").Value = 71598
$ws.Range("J125").Value = 71598
$ws.Range("L125").Value = 71598
$ws.Range("N125").Value = -81438
$ws.Range("H136").Value = 3232.6667
$ws.Range("I136").Value = 1199.5
$ws.Range("J136").Value = 4249.25
$ws.Range("K136").Value = 3598.5
$ws.Range("L136").Value = 12747.75
$ws.Range("M136").Value = -1048.5
$ws.Range("N136").Value = -17847.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 35000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 35000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 35000
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = -35226
$ws.Range("H17").Value = 9000
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = ""
$ws.Range("H132").Value = 580.7
$ws.Range("I132").Value = 467.44446
$ws.Range("K132").Value = 1402.33338
$ws.Range("M132").Value = 1127.66662
$ws.Range("H136").Value = 3126.6667
$ws.Range("I136").Value = 3196.3076
$ws.Range("J136").Value = 2945.6
$ws.Range("K136").Value = 9588.9228
$ws.Range("L136").Value = 8836.799999999999
$ws.Range("M136").Value = -7038.9228
$ws.Range("N136").Value = -13936.8
